$d = $word.ActiveDocument
$s = $d.Styles.Item("Heading1")
$s.ParagraphFormat.WidowControl = 0
Write-Output $s.ParagraphFormat.WidowControl
